$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 12141.111
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 20254
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 60762
$ws.Range("M17").Value = -5832
$ws.Range("N17").Value = -61098

$ws.Range("H19").Value = 747.6923
$ws.Range("I19").Value = 466.33334
$ws.Range("J19").Value = 832.1
$ws.Range("K19").Value = 466.33334
$ws.Range("L19").Value = 832.1
$ws.Range("M19").Value = -291.33334
$ws.Range("N19").Value = -1182.1

$ws.Range("J54").Value = 80000
$ws.Range("K54").Value = 80000
$ws.Range("L54").Value = 80000
$ws.Range("M54").Value = -79514
$ws.Range("N54").Value = -80972

$ws.Range("H76").Value = 490
$ws.Range("I76").Value = 490
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 490
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -175

$ws.Range("H79").Value = 490
$ws.Range("I79").Value = 490
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 490
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 602

$ws.Range("H132").Value = 1218.3334
$ws.Range("I132").Value = 712
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 2136
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = 394

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1827
$ws.Range("I26").Value = 463.5
$ws.Range("J26").Value = 10008
$ws.Range("K26").Value = 463.5
$ws.Range("L26").Value = 10008
$ws.Range("M26").Value = -133.5
$ws.Range("N26").Value = -10668

$ws.Range("H32").Value = 199.9
$ws.Range("I32").Value = 199.9
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 199.9
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 87.09999999999999

$ws.Range("H97").Value = 1934.6666
$ws.Range("I97").Value = 1934.6666
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1934.6666
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1438.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4274.75
$ws.Range("I20").Value = 3599.75
$ws.Range("J20").Value = 4949.75
$ws.Range("K20").Value = 3599.75
$ws.Range("L20").Value = 4949.75
$ws.Range("M20").Value = -3352.75
$ws.Range("N20").Value = -5443.75

$ws.Range("H99").Value = 5066.1113
$ws.Range("I99").Value = 5599.8335
$ws.Range("J99").Value = 3998.6667
$ws.Range("K99").Value = 5599.8335
$ws.Range("L99").Value = 3998.6667
$ws.Range("M99").Value = -4101.8335

$ws.Range("H105").Value = 2076.5
$ws.Range("I105").Value = 1995
$ws.Range("J105").Value = 2103.6667
$ws.Range("K105").Value = 1995
$ws.Range("L105").Value = 2103.6667
$ws.Range("M105").Value = -248
$ws.Range("N105").Value = -5597.6667

$ws.Range("H134").Value = 5432.5
$ws.Range("I134").Value = 5119
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 15357
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -12822
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 6200.5557
$ws.Range("I39").Value = 3674.6667
$ws.Range("J39").Value = 11252.333
$ws.Range("K39").Value = 3674.6667
$ws.Range("L39").Value = 11252.333
$ws.Range("M39").Value = -3283.6667

$ws.Range("H49").Value = 6200.5557
$ws.Range("I49").Value = 3674.6667
$ws.Range("J49").Value = 11252.333
$ws.Range("K49").Value = 3674.6667
$ws.Range("L49").Value = 11252.333
$ws.Range("M49").Value = -3492.6667

$ws.Range("H132").Value = 4342
$ws.Range("I132").Value = 2012
$ws.Range("J132").Value = 5507
$ws.Range("K132").Value = 6036
$ws.Range("L132").Value = 16521
$ws.Range("M132").Value = -3506
$ws.Range("N132").Value = -21581

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 32319.857
$ws.Range("I32").Value = 4620.5
$ws.Range("J32").Value = 43399.6
$ws.Range("K32").Value = 13861.5
$ws.Range("L32").Value = 130198.8
$ws.Range("M32").Value = -13578.5

$ws.Range("H34").Value = 1549.6666
$ws.Range("I34").Value = 825
$ws.Range("J34").Value = 2999
$ws.Range("K34").Value = 2475
$ws.Range("L34").Value = 8997
$ws.Range("M34").Value = -2391
$ws.Range("N34").Value = -9165

$ws.Range("H70").Value = 6816.1665
$ws.Range("I70").Value = 6816.1665
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 20448.4995
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -20133.4995

$ws.Range("H73").Value = 6816.1665
$ws.Range("I73").Value = 6816.1665
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 20448.4995
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -19356.4995

$ws.Range("H122").Value = 656.4286
$ws.Range("I122").Value = 479
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 4311
$ws.Range("L122").Value = 9900
$ws.Range("M122").Value = -1861
$ws.Range("N122").Value = -14800

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H23").Value = 963.75
$ws.Range("I23").Value = 830
$ws.Range("J23").Value = 1900
$ws.Range("K23").Value = 830
$ws.Range("L23").Value = 1900
$ws.Range("M23").Value = -607
$ws.Range("N23").Value = -2346

$ws.Range("H97").Value = 1832.6666
$ws.Range("I97").Value = 1799
$ws.Range("J97").Value = 1900
$ws.Range("K97").Value = 1799
$ws.Range("L97").Value = 1900
$ws.Range("M97").Value = -1303
$ws.Range("N97").Value = -2892

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 797.7143
$ws.Range("I22").Value = 714.8333
$ws.Range("J22").Value = 859.875
$ws.Range("K22").Value = 714.8333
$ws.Range("L22").Value = 859.875
$ws.Range("M22").Value = -419.8333
$ws.Range("N22").Value = -1449.875

$ws.Range("H27").Value = 797.7143
$ws.Range("I27").Value = 714.8333
$ws.Range("J27").Value = 859.875
$ws.Range("K27").Value = 714.8333
$ws.Range("L27").Value = 859.875
$ws.Range("M27").Value = -607.8333
$ws.Range("N27").Value = -1073.875

$ws.Range("H61").Value = 18497.834
$ws.Range("I61").Value = 1995.5
$ws.Range("J61").Value = 51502.5
$ws.Range("K61").Value = 1995.5
$ws.Range("L61").Value = 51502.5
$ws.Range("M61").Value = -1793.5
$ws.Range("N61").Value = -51906.5

$ws.Range("H82").Value = 911
$ws.Range("I82").Value = 822.1
$ws.Range("J82").Value = 1800
$ws.Range("K82").Value = 822.1
$ws.Range("L82").Value = 1800
$ws.Range("M82").Value = -461.1

$ws.Range("H85").Value = 911
$ws.Range("I85").Value = 822.1
$ws.Range("J85").Value = 1800
$ws.Range("K85").Value = 822.1
$ws.Range("L85").Value = 1800
$ws.Range("M85").Value = 425.9

$ws.Range("H113").Value = 18497.834
$ws.Range("I113").Value = 1995.5
$ws.Range("J113").Value = 51502.5
$ws.Range("K113").Value = 1995.5
$ws.Range("L113").Value = 51502.5
$ws.Range("M113").Value = 174.5
$ws.Range("N113").Value = -55842.5

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H132").Value = 2674.75
$ws.Range("I132").Value = 2674.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8024.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5494.25
